# Add section "2.3 多级反馈队列" after the time-slice-rotation paragraph,
# right before the trailing (bookmark) paragraph at the end of the document.
#
# Each inner array is one new paragraph; its elements are the individual
# <w:r> run texts that paragraph should end up with (an empty array means
# an empty paragraph). This mirrors the target OOXML exactly, including
# the run boundaries.
$paraRuns = @(
    ,@()
    ,@('2.3 ', '多级反馈队列')
    ,@('如果一个进程需要执行', ' 100 ', '个时间片，如果采用轮转调度算法，那么需要交换', ' 100 ', '次。多级队列是为这种需要连续执行多个时间片的进程考虑，它设置了多个队列，每个队列时间片大小都不同，例如', ' 1,2,4,8,..', '。进程在第一个队列没执行完，就会被移到下一个队列。这种方式下，之前的进程只需要', ' 7 ', '（包括最初的装入）的交换。')
    ,@('每个队列的优先权也不同，最上面的优先权最高。因此只有上一个队列没有进程在排队，才能调度当前队列上的进程。')
    ,@()
)

$d = $word.ActiveDocument

# Locate the paragraph that ends the "时间片轮转" discussion -- the new
# content is inserted right after it.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*在进程切换上就会花过多时间*") {
        $target = $p
    }
}

$insertStart = $target.Range.End

# Build the whole block (all new paragraphs) as a single string, joined by
# carriage returns, and drop it in with one Range.Text assignment. Doing it
# this way (rather than many small InsertParagraphAfter/InsertAfter calls)
# keeps truly-empty paragraphs free of a stray placeholder run, and keeps
# the runs of a paragraph merged into one contiguous string for now -- we
# fix the run boundaries below.
$blockParts = @()
foreach ($runs in $paraRuns) {
    $blockParts += ($runs -join '')
}
$block = [char]13 + ($blockParts -join [char]13)

$r = $d.Range($insertStart, $insertStart)
$r.Text = $block

# Now split multi-run paragraphs back into their individual runs. Selecting
# a non-empty sub-range and toggling a character property (Bold on, then
# back off) forces the engine to materialize that sub-range as its own
# run without leaving any residual formatting behind.
$pos = $insertStart + 1
foreach ($runs in $paraRuns) {
    $paraLen = 0
    foreach ($t in $runs) { $paraLen += $t.Length }

    if ($runs.Count -gt 1) {
        $cum = 0
        $boundaries = @()
        for ($i = 0; $i -lt $runs.Count - 1; $i++) {
            $cum = $cum + $runs[$i].Length
            $boundaries += $cum
        }
        for ($i = $boundaries.Count - 1; $i -ge 0; $i--) {
            $b = $boundaries[$i]
            $sub = $d.Range($pos + $b, $pos + $paraLen)
            $sub.Bold = 1
            $sub.Bold = 0
        }
    }

    $pos = $pos + $paraLen + 1
}
